# Updates the "cryptos" price/volume table to the latest scrape values.
#
# Columns: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Rows 42/43 (Stacks <-> OKB) also swap position in addition to their values
# changing, matching the upstream ranking re-sort.
#
# Price strings that look like plain numbers (e.g. "0.998", "26.50") are
# written with a leading apostrophe so Excel keeps them as literal text
# instead of silently re-typing them as numeric cells (which would also
# collapse formatting such as trailing zeros, e.g. "26.50" -> 26.5). The
# apostrophe itself is not stored as part of the cell's value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.468.94"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.456.34"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'580.03"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").Value = "'143.84"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.84%  "
$ws.Range("D9").Value = "2.452.66"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +2.62%  "
$ws.Range("D12").Value = "'5.22"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "'0.348"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "'26.50"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "'0.0000178"
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").Value = "2.897.94"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "62.152.69"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "2.452.66"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "'10.90"
$ws.Range("E19").Value = "  -1.85%  "
$ws.Range("D20").Value = "'7.14"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "'328.81"
$ws.Range("E21").Value = "  +1.49%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  -3.74%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "'65.72"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("D26").Value = "'9.36"
$ws.Range("E26").Value = "  +5.34%  "
$ws.Range("D27").Value = "'590.75"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("D28").Value = "0.0₃0970"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'8.07"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "'1.88"
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("D34").Value = "'0.135"
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").Value = "'4.95"
$ws.Range("E35").Value = "  -1.23%  "
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("D37").Value = "'1.46"
$ws.Range("E37").Value = "  -1.04%  "
$ws.Range("D38").Value = "'0.380"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").Value = "'155.04"
$ws.Range("E39").Value = "  +5.19%  "
$ws.Range("D40").Value = "'5.35"
$ws.Range("E40").Value = "  +2.31%  "
$ws.Range("D41").Value = "'18.45"
$ws.Range("E41").Value = "  -1.24%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").Value = "'43.13"
$ws.Range("E42").Value = "  +2.49%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.73"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D45").Value = "'2.47"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "'143.36"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").Value = "'3.67"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").Value = "0.0₆0261"
$ws.Range("E48").Value = "  +20.17%  "
$ws.Range("D49").Value = "'0.610"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").Value = "'0.0525"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").Value = "'19.94"
$ws.Range("E51").Value = "  -1.50%  "
